$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F2").Value = 1664.71
$wsSummary.Range("A3").Value = 75.1
$wsSummary.Range("E3").Value = 75.1
$wsSummary.Range("F3").Value = 29.17
# Convert cells from 0.00 number format to General (matches existing style 14)
$wsSummary.Range("A3").NumberFormat = "general"
$wsSummary.Range("E3").NumberFormat = "general"
$wsSummary.Range("F3").NumberFormat = "general"

# ---- Repayment schedule sheet ----
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("F3").Value = 834.44
$wsRepay.Range("G3").Value = 4165.56
$wsRepay.Range("H3").Value = 12.5
$wsRepay.Range("K3").Value = 846.94
$wsRepay.Range("P3").Value = 846.94
$wsRepay.Range("F4").Value = 830.27
$wsRepay.Range("G4").Value = 3335.29
$wsRepay.Range("H4").Value = 16.67
$wsRepay.Range("K4").Value = 846.94
$wsRepay.Range("P4").Value = 846.94
$wsRepay.Range("F5").Value = 826.11
$wsRepay.Range("G5").Value = 2509.18
$wsRepay.Range("H5").Value = 20.83
$wsRepay.Range("K5").Value = 846.94
$wsRepay.Range("P5").Value = 846.94
$wsRepay.Range("F6").Value = 834.39
$wsRepay.Range("G6").Value = 1674.79
$wsRepay.Range("H6").Value = 12.55
$wsRepay.Range("K6").Value = 846.94
$wsRepay.Range("P6").Value = 846.94
$wsRepay.Range("F7").Value = 838.57
$wsRepay.Range("G7").Value = 836.22
$wsRepay.Range("H7").Value = 8.37
$wsRepay.Range("K7").Value = 846.94
$wsRepay.Range("P7").Value = 846.94
$wsRepay.Range("F8").Value = 836.22
$wsRepay.Range("G8").Value = 0
$wsRepay.Range("H8").Value = 4.18
$wsRepay.Range("K8").Value = 840.4
$wsRepay.Range("P8").Value = 840.4
# Convert cells from 0.00 / #,##0.00 number formats to General (matches existing style 14)
$wsRepay.Range("F3").NumberFormat = "general"
$wsRepay.Range("K3").NumberFormat = "general"
$wsRepay.Range("P3").NumberFormat = "general"
$wsRepay.Range("F4").NumberFormat = "general"
$wsRepay.Range("K4").NumberFormat = "general"
$wsRepay.Range("P4").NumberFormat = "general"
$wsRepay.Range("F5").NumberFormat = "general"
$wsRepay.Range("K5").NumberFormat = "general"
$wsRepay.Range("P5").NumberFormat = "general"
$wsRepay.Range("F6").NumberFormat = "general"
$wsRepay.Range("K6").NumberFormat = "general"
$wsRepay.Range("P6").NumberFormat = "general"
$wsRepay.Range("F7").NumberFormat = "general"
$wsRepay.Range("K7").NumberFormat = "general"
$wsRepay.Range("P7").NumberFormat = "general"
$wsRepay.Range("F8").NumberFormat = "general"
$wsRepay.Range("K8").NumberFormat = "general"
$wsRepay.Range("P8").NumberFormat = "general"
$wsRepay.Range("H3").NumberFormat = "general"
$wsRepay.Range("H4").NumberFormat = "general"
$wsRepay.Range("H5").NumberFormat = "general"
$wsRepay.Range("H6").NumberFormat = "general"
$wsRepay.Range("H7").NumberFormat = "general"
$wsRepay.Range("H8").NumberFormat = "general"
$wsRepay.Range("G7").NumberFormat = "general"
$wsRepay.Range("G8").NumberFormat = "general"

# Remove now-unused rows 9:18 (dimension shrinks to A1:P8)
$wsRepay.Rows("9:18").Delete()
# Re-apply explicit row height on row 8 (matches ht="15" customHeight="1")
$wsRepay.Rows("8:8").RowHeight = 15

# ---- Selections / active sheet sequence (order matters for tabSelected) ----
$wsSummary.Activate()
$wsSummary.Range("A7:XFD14").Select()

$wsRepay.Activate()
$wsRepay.Range("A9:XFD13").Select()

$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()

